$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("relays")

# Update the Max. Slip Voltage [%] (dV) values from 5 to 10 for rows 2-7 (column C)
$ws.Range("C2:C7").Value = 10

# Update the active selection to E12 as per the recorded workbook state
$ws.Range("E12").Select()
